# Update cryptos list data (prices, volume %, and a few coin row swaps)
# as scraped on Sat Jul 22 11:52:07 UTC 2023

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.895.92"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.889.17"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7750"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.96"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3116"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.63"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07174"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08605"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +6.39%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.970.20"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.55%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7635"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.378"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.79"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.178"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.946.39"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.78"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.50"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007821"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.194.94"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9984"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.002"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.99%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1640"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.384"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.06"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.042"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.444"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.541"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.522"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.101"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05433"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.240"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7461"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.66%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.695"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.13%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.97%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4472"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.116.23"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.40%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "73.19"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.087"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.05%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.47"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.866"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.644"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.086.93"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.982"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.77%  "
